# Fruta / hortaliza, semanal
# A new weekly observation is inserted as row 11 (shifting the existing
# rows 11-13 down to rows 12-14), and the new row 11 is populated with
# the latest weekly price data for "Mora" at Vega Central Mapocho de
# Santiago / Provincia de Curicó.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing existing rows 11:13 down to 12:14.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44586
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100101008
$ws.Range("J11").Value = "Mora"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 5000
$ws.Range("O11").Value = 5000
$ws.Range("P11").Value = 5000
$ws.Range("Q11").Value = "$/bandeja 2 kilos"
$ws.Range("R11").Value = "Provincia de Curicó"
$ws.Range("S11").Value = 2500
$ws.Range("T11").Value = 2
